# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh described in the commit diff:
# updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ cells (columns H-N) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, matching the public diff exactly.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 1999.6666
$ws.Range("J49").Value = 1999.5
$ws.Range("L49").Value = 5998.5
$ws.Range("N49").Value = -6270.5
$ws.Range("H64").Value = 4888
$ws.Range("J64").Value = 4999.6665
$ws.Range("L64").Value = 4999.6665
$ws.Range("N64").Value = -5495.6665
$ws.Range("H67").Value = 4888
$ws.Range("J67").Value = 4999.6665
$ws.Range("L67").Value = 4999.6665
$ws.Range("N67").Value = -6715.6665
$ws.Range("H107").Value = 749
$ws.Range("I107").Value = 749
$ws.Range("K107").Value = 749
$ws.Range("M107").Value = 1171
$ws.Range("H132").Value = 2525.463
$ws.Range("I132").Value = 2559.68
$ws.Range("K132").Value = 7679.039999999999
$ws.Range("M132").Value = -5149.039999999999
$ws.Range("H138").Value = 2495.0454
$ws.Range("J138").Value = 4580.3887
$ws.Range("L138").Value = 13741.1661
$ws.Range("N138").Value = -24021.1661

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2006.8
$ws.Range("J45").Value = 2675.6667
$ws.Range("L45").Value = 2675.6667
$ws.Range("N45").Value = -3429.6667
$ws.Range("H61").Value = 37039104
$ws.Range("I61").Value = 47620332
$ws.Range("K61").Value = 47620332
$ws.Range("M61").Value = -47620120
$ws.Range("H63").Value = 6187.154
$ws.Range("I63").Value = 6281.6665
$ws.Range("J63").Value = 5974.5
$ws.Range("K63").Value = 6281.6665
$ws.Range("L63").Value = 5974.5
$ws.Range("M63").Value = -5595.6665
$ws.Range("N63").Value = -7346.5
$ws.Range("H66").Value = 6187.154
$ws.Range("I66").Value = 6281.6665
$ws.Range("J66").Value = 5974.5
$ws.Range("K66").Value = 31408.3325
$ws.Range("L66").Value = 29872.5
$ws.Range("M66").Value = -27976.3325
$ws.Range("N66").Value = -36736.5
$ws.Range("H122").Value = 17546348
$ws.Range("I122").Value = 27779222
$ws.Range("J122").Value = 4278.5713
$ws.Range("K122").Value = 83337666
$ws.Range("L122").Value = 12835.7139
$ws.Range("M122").Value = -83335216
$ws.Range("N122").Value = -17735.7139
$ws.Range("H132").Value = 2631
$ws.Range("I132").Value = 2631
$ws.Range("K132").Value = 7893
$ws.Range("M132").Value = -5363
$ws.Range("H136").Value = 37039104
$ws.Range("I136").Value = 47620332
$ws.Range("K136").Value = 142860996
$ws.Range("M136").Value = -142858446

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 41917.445
$ws.Range("I82").Value = 15564.25
$ws.Range("J82").Value = 63000
$ws.Range("K82").Value = 15564.25
$ws.Range("L82").Value = 63000
$ws.Range("M82").Value = -15181.25
$ws.Range("N82").Value = -63766
$ws.Range("H85").Value = 41917.445
$ws.Range("I85").Value = 15564.25
$ws.Range("J85").Value = 63000
$ws.Range("K85").Value = 15564.25
$ws.Range("L85").Value = 63000
$ws.Range("M85").Value = -14238.25
$ws.Range("N85").Value = -65652
$ws.Range("H97").Value = 3181.5
$ws.Range("I97").Value = 3181.5
$ws.Range("K97").Value = 3181.5
$ws.Range("M97").Value = -2190.5
$ws.Range("H99").Value = 1394.8889
$ws.Range("I99").Value = 1017.3333
$ws.Range("K99").Value = 1017.3333
$ws.Range("M99").Value = 480.6667
$ws.Range("H105").Value = 2301
$ws.Range("I105").Value = 2181.8
$ws.Range("J105").Value = 2499.6667
$ws.Range("K105").Value = 2181.8
$ws.Range("L105").Value = 2499.6667
$ws.Range("M105").Value = -434.8000000000002
$ws.Range("N105").Value = -5993.6667
$ws.Range("H126").Value = 49999
$ws.Range("J126").Value = 49999
$ws.Range("L126").Value = 49999
$ws.Range("N126").Value = -59879
$ws.Range("H134").Value = 1958.4
$ws.Range("I134").Value = 1958.4
$ws.Range("K134").Value = 5875.200000000001
$ws.Range("M134").Value = -3340.200000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1503.5
$ws.Range("I25").Value = 338
$ws.Range("K25").Value = 338
$ws.Range("M25").Value = -164
$ws.Range("H39").Value = 370
$ws.Range("I39").Value = 370
$ws.Range("K39").Value = 370
$ws.Range("M39").Value = 21
$ws.Range("H49").Value = 370
$ws.Range("I49").Value = 370
$ws.Range("K49").Value = 370
$ws.Range("M49").Value = -188
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("H99").Value = 2622.7334
$ws.Range("I99").Value = 2266.5557
$ws.Range("J99").Value = 3157
$ws.Range("K99").Value = 2266.5557
$ws.Range("L99").Value = 3157
$ws.Range("M99").Value = -768.5556999999999
$ws.Range("N99").Value = -6153
$ws.Range("H105").Value = 1511.8667
$ws.Range("I105").Value = 845.8
$ws.Range("K105").Value = 845.8
$ws.Range("M105").Value = 901.2
$ws.Range("H122").Value = 2808.077
$ws.Range("I122").Value = 2001
$ws.Range("J122").Value = 5498.3335
$ws.Range("K122").Value = 6003
$ws.Range("L122").Value = 16495.0005
$ws.Range("M122").Value = -3553
$ws.Range("N122").Value = -21395.0005
$ws.Range("H126").Value = 2622.7334
$ws.Range("I126").Value = 2266.5557
$ws.Range("J126").Value = 3157
$ws.Range("K126").Value = 6799.6671
$ws.Range("L126").Value = 9471
$ws.Range("M126").Value = -4329.6671
$ws.Range("N126").Value = -14411
$ws.Range("N97").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 40000
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 120000
$ws.Range("N37").Value = -120224
$ws.Range("H55").Value = 897.5925999999999
$ws.Range("J55").Value = 2250
$ws.Range("L55").Value = 6750
$ws.Range("N55").Value = -7104

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 3
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 166
$ws.Range("H58").Value = 21989.6
$ws.Range("I58").Value = 17500
$ws.Range("K58").Value = 17500
$ws.Range("M58").Value = -17223
$ws.Range("H122").Value = 2920
$ws.Range("I122").Value = 2587.7
$ws.Range("K122").Value = 7763.099999999999
$ws.Range("M122").Value = -5313.099999999999
$ws.Range("H126").Value = 9644.214
$ws.Range("I126").Value = 11183.546
$ws.Range("K126").Value = 33550.638
$ws.Range("M126").Value = -31080.638

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20835890
$ws.Range("I7").Value = 26317796
$ws.Range("J7").Value = 4650
$ws.Range("K7").Value = 26317796
$ws.Range("L7").Value = 4650
$ws.Range("M7").Value = -26317684
$ws.Range("N7").Value = -4874
$ws.Range("H40").Value = 4998.8335
$ws.Range("I40").Value = 4996.5
$ws.Range("K40").Value = 4996.5
$ws.Range("M40").Value = -4860.5
$ws.Range("H55").Value = 558.7059
$ws.Range("I55").Value = 573.75
$ws.Range("J55").Value = 554.0769
$ws.Range("K55").Value = 573.75
$ws.Range("L55").Value = 554.0769
$ws.Range("M55").Value = -400.75
$ws.Range("N55").Value = -900.0769
$ws.Range("H61").Value = 1398.8
$ws.Range("I61").Value = 1398.8
$ws.Range("K61").Value = 1398.8
$ws.Range("M61").Value = -1196.8
$ws.Range("H100").Value = 3471.2856
$ws.Range("I100").Value = 3459.8
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 3459.8
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -2918.8
$ws.Range("N100").Value = -4582
$ws.Range("H113").Value = 1398.8
$ws.Range("I113").Value = 1398.8
$ws.Range("K113").Value = 1398.8
$ws.Range("M113").Value = 771.2
$ws.Range("H122").Value = 5666.6665
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 6285.7144
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 18857.1432
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -23757.1432
$ws.Range("H126").Value = 20835890
$ws.Range("I126").Value = 26317796
$ws.Range("J126").Value = 4650
$ws.Range("K126").Value = 78953388
$ws.Range("L126").Value = 13950
$ws.Range("M126").Value = -78950918
$ws.Range("N126").Value = -18890

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 211.33333
$ws.Range("I7").Value = 211.33333
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 211.33333
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -98.33332999999999
$ws.Range("H21").Value = 39013.6
$ws.Range("I21").Value = 35000
$ws.Range("J21").Value = 40017
$ws.Range("K21").Value = 35000
$ws.Range("L21").Value = 40017
$ws.Range("M21").Value = -34765
$ws.Range("N21").Value = -40487
$ws.Range("H35").Value = 39013.6
$ws.Range("I35").Value = 35000
$ws.Range("J35").Value = 40017
$ws.Range("K35").Value = 35000
$ws.Range("L35").Value = 40017
$ws.Range("M35").Value = -34710
$ws.Range("N35").Value = -40597
$ws.Range("H122").Value = 2224.75
$ws.Range("I122").Value = 1950
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400
$ws.Range("H126").Value = 1694.4546
$ws.Range("I126").Value = 1613.9
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4841.700000000001
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2371.700000000001
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 3796.6667
$ws.Range("I132").Value = 3577.0386
$ws.Range("J132").Value = 5224.25
$ws.Range("K132").Value = 10731.1158
$ws.Range("L132").Value = 15672.75
$ws.Range("M132").Value = -8201.1158
$ws.Range("N132").Value = -20732.75
$ws.Range("N7").ClearContents()

